$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 269 (shifts existing rows 269:283 down to 270:284)
$ws.Rows.Item(269).Insert()

# Populate the newly inserted row 269 with the new weekly record
$ws.Cells.Item(269, 1).Value = 7
$ws.Cells.Item(269, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(269, 3).Value = "Ñuble"
$ws.Cells.Item(269, 4).Value = 44516
$ws.Cells.Item(269, 5).Value = 16
$ws.Cells.Item(269, 6).Value = 100114014
$ws.Cells.Item(269, 7).Value = "Betarraga"
$ws.Cells.Item(269, 8).Value = "Sin especificar"
$ws.Cells.Item(269, 9).Value = "Primera"
$ws.Cells.Item(269, 10).Value = 360
$ws.Cells.Item(269, 11).Value = 700
$ws.Cells.Item(269, 12).Value = 800
$ws.Cells.Item(269, 13).Value = 750
$ws.Cells.Item(269, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(269, 15).Value = "Región del Maule"
$ws.Cells.Item(269, 16).Value = 150
$ws.Cells.Item(269, 17).Value = 5
$ws.Cells.Item(269, 18).Value = "Hortaliza"
